$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the URL value cell (C7) from the old placeholder text to the actual demo store URL
$ws.Range("C7").Value = "https://demo.nopcommerce.com/"

# Update the selected/active cell in the sheet view
$ws.Range("F8").Select()
